$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EDM DATA")

$data = @(
    @("2024-12-19 00:23:26", -0.1208256299179629, -0.001662890008911996, 0.008036789312443157),
    @("2024-12-19 00:23:27", -0.1211997625412187, -0.001651754980391996, 0.008007692455991408),
    @("2024-12-19 00:23:28", -0.1199895248381654, -0.001755968724585996, 0.008427914115750119),
    @("2024-12-19 00:23:29", -0.1204580039490247, -0.001647857720409997, 0.00793990607170313),
    @("2024-12-19 00:23:30", -0.1185385409253649, -0.001844694656383996, 0.008746696520823053),
    @("2024-12-19 00:23:31", -0.120441737313231,  -0.001744783082299996, 0.008405788226677825),
    @("2024-12-19 00:23:32", -0.1192868061718764, -0.001651552525327996, 0.007880337038858943),
    @("2024-12-19 00:23:33", -0.1192477662459714, -0.001751109803049995, 0.008352637298605397),
    @("2024-12-19 00:23:34", -0.1186361407401273, -0.001675138540283995, 0.007949278864973737)
)

$startRow = 140
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
